# Remove `event_name` from required columns
#
# The "column_names" sheet holds a lookup table (Table7) mapping the raw
# source-system column names (column A) to the cleaned/renamed columns used
# downstream (column B). Row 6 mapped "EventName" -> "event"; the author
# removed that mapping entirely, so the table (and its backing range) shrink
# from 12 to 11 rows and every row below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("column_names")

# Deleting the sheet row also contracts the ListObject/table range
# (A1:B12 -> A1:B11) and re-numbers the shared strings automatically.
$ws.Rows.Item(6).Delete()

# The author had this sheet active (and selection reset to the top-left
# cell) when the workbook was last saved, moving the active tab away from
# "settings" back to "column_names".
$ws.Activate()
$ws.Range("A1").Select()
